$p = $ppt.ActivePresentation

# Use the same layout as the existing content slides (title + content)
$layout = $p.Slides.Item($p.Slides.Count).CustomLayout

# Append a new slide at the end of the deck
$newSlide = $p.Slides.Add($p.Slides.Count + 1, $layout)

# Title placeholder -> "END"
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "END"
